# Scheduled refresh of market-price-derived profit figures (columns H-N)
# across the per-job profit sheets. Values below come from the upstream
# market-board snapshot used to recompute currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) for the affected Leve rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 780.3077
$ws.Range("I19").Value = 518.8570999999999
$ws.Range("J19").Value = 1085.3334
$ws.Range("K19").Value = 518.8570999999999
$ws.Range("L19").Value = 1085.3334
$ws.Range("M19").Value = -343.8570999999999
$ws.Range("N19").Value = -1435.3334

$ws.Range("H28").Value = 1047.2778
$ws.Range("I28").Value = 1203.8667
$ws.Range("J28").Value = 264.33334
$ws.Range("K28").Value = 1203.8667
$ws.Range("L28").Value = 264.33334
$ws.Range("M28").Value = -718.8667
$ws.Range("N28").Value = -1234.33334

$ws.Range("H100").Value = 1404.2593
$ws.Range("I100").Value = 1326.0714
$ws.Range("J100").Value = 1488.4615
$ws.Range("K100").Value = 1326.0714
$ws.Range("L100").Value = 1488.4615
$ws.Range("M100").Value = -785.0714
$ws.Range("N100").Value = -2570.4615

$ws.Range("H115").Value = 1077.2858
$ws.Range("I115").Value = 590.1667
$ws.Range("J115").Value = 4000
$ws.Range("K115").Value = 1770.5001
$ws.Range("L115").Value = 12000
$ws.Range("M115").Value = -203.5001
$ws.Range("N115").Value = -15134

$ws.Range("H118").Value = 445
$ws.Range("I118").Value = 292.5
$ws.Range("J118").Value = 750
$ws.Range("K118").Value = 877.5
$ws.Range("L118").Value = 2250
$ws.Range("M118").Value = 779.5
$ws.Range("N118").Value = -5564

$ws.Range("H139").Value = 37797.145
$ws.Range("J139").Value = 37797.145
$ws.Range("L139").Value = 37797.145
$ws.Range("N139").Value = -48077.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3604.4
$ws.Range("I2").Value = 2893.7778
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 2893.7778
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -2780.7778
$ws.Range("N2").Value = -10226

$ws.Range("H32").Value = 5425.939
$ws.Range("I32").Value = 4689.1953
$ws.Range("K32").Value = 4689.1953
$ws.Range("M32").Value = -4402.1953

$ws.Range("H116").Value = 3604.4
$ws.Range("I116").Value = 2893.7778
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 2893.7778
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -599.7777999999998
$ws.Range("N116").Value = -14588

$ws.Range("H122").Value = 1781.5
$ws.Range("I122").Value = 1775.6666
$ws.Range("K122").Value = 5326.9998
$ws.Range("M122").Value = -2876.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3604.4
$ws.Range("I3").Value = 2893.7778
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 2893.7778
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -2779.7778
$ws.Range("N3").Value = -10228

$ws.Range("H107").Value = 10478.786
$ws.Range("I107").Value = 1119
$ws.Range("K107").Value = 1119
$ws.Range("M107").Value = 801

$ws.Range("H135").Value = 45500
$ws.Range("J135").Value = 45500
$ws.Range("L135").Value = 45500
$ws.Range("N135").Value = -55640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1316.7894
$ws.Range("I16").Value = 1270.9286
$ws.Range("K16").Value = 1270.9286
$ws.Range("M16").Value = -983.9286

$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -2609
$ws.Range("N39").ClearContents()

$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 3000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 3000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -2818
$ws.Range("N49").ClearContents()

$ws.Range("H113").Value = 1316.7894
$ws.Range("I113").Value = 1270.9286
$ws.Range("K113").Value = 1270.9286
$ws.Range("M113").Value = 899.0714

$ws.Range("H134").Value = 3050.389
$ws.Range("I134").Value = 2188.111
$ws.Range("J134").Value = 3912.6667
$ws.Range("K134").Value = 6564.333
$ws.Range("L134").Value = 11738.0001
$ws.Range("M134").Value = -4029.333
$ws.Range("N134").Value = -16808.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 465.44446
$ws.Range("I121").Value = 384.14285
$ws.Range("J121").Value = 750
$ws.Range("K121").Value = 1152.42855
$ws.Range("L121").Value = 2250
$ws.Range("M121").Value = 157.5714499999999
$ws.Range("N121").Value = -4870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1089871.1
$ws.Range("I107").Value = 1838768.5
$ws.Range("J107").Value = 565.8182
$ws.Range("K107").Value = 1838768.5
$ws.Range("L107").Value = 565.8182
$ws.Range("M107").Value = -1836848.5
$ws.Range("N107").Value = -4405.8182

$ws.Range("H123").Value = 25546.223
$ws.Range("J123").Value = 25546.223
$ws.Range("L123").Value = 25546.223
$ws.Range("N123").Value = -30446.223

$ws.Range("H132").Value = 2404.261
$ws.Range("I132").Value = 1971.4667
$ws.Range("J132").Value = 3215.75
$ws.Range("K132").Value = 5914.4001
$ws.Range("L132").Value = 9647.25
$ws.Range("M132").Value = -3384.4001
$ws.Range("N132").Value = -14707.25

$ws.Range("H139").Value = 18950
$ws.Range("J139").Value = 18950
$ws.Range("L139").Value = 18950
$ws.Range("N139").Value = -29230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1610.2
$ws.Range("I61").Value = 1128.8572
$ws.Range("K61").Value = 1128.8572
$ws.Range("M61").Value = -926.8571999999999

$ws.Range("H113").Value = 1610.2
$ws.Range("I113").Value = 1128.8572
$ws.Range("K113").Value = 1128.8572
$ws.Range("M113").Value = 1041.1428

$ws.Range("H132").Value = 2833.8235
$ws.Range("I132").Value = 2740.4614
$ws.Range("J132").Value = 3137.25
$ws.Range("K132").Value = 8221.3842
$ws.Range("L132").Value = 9411.75
$ws.Range("M132").Value = -5691.3842
$ws.Range("N132").Value = -14471.75

$ws.Range("H134").Value = 19976.334
$ws.Range("J134").Value = 19976.334
$ws.Range("L134").Value = 19976.334
$ws.Range("N134").Value = -30116.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 956.8182
$ws.Range("I107").Value = 978.125
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 2934.375
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -1014.375
$ws.Range("N107").Value = -6540

$ws.Range("H113").Value = 264.64285
$ws.Range("I113").Value = 262.2
$ws.Range("J113").Value = 270.75
$ws.Range("K113").Value = 786.5999999999999
$ws.Range("L113").Value = 812.25
$ws.Range("M113").Value = 1383.4
$ws.Range("N113").Value = -5152.25

$ws.Range("H138").Value = 43933.332
$ws.Range("J138").Value = 43933.332
$ws.Range("L138").Value = 43933.332
$ws.Range("N138").Value = -54213.332
